function Set-TextValue {
    param($range, $value)
    $orig = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $orig
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-TextValue $ws.Range('D2') '70.416.94'
Set-TextValue $ws.Range('E2') '  -2.58%  '

Set-TextValue $ws.Range('D3') '3.625.42'
Set-TextValue $ws.Range('E3') '  +0.44%  '

Set-TextValue $ws.Range('E4') '  +0.14%  '

Set-TextValue $ws.Range('D5') '579.62'
Set-TextValue $ws.Range('E5') '  -2.92%  '

Set-TextValue $ws.Range('D6') '174.97'
Set-TextValue $ws.Range('E6') '  -4.81%  '

Set-TextValue $ws.Range('D7') '3.618.19'
Set-TextValue $ws.Range('E7') '  +0.51%  '

Set-TextValue $ws.Range('D8') '0.608'
Set-TextValue $ws.Range('E8') '  -0.09%  '

Set-TextValue $ws.Range('E9') '  -0.03%  '

Set-TextValue $ws.Range('D10') '0.195'
Set-TextValue $ws.Range('E10') '  -5.68%  '

Set-TextValue $ws.Range('D11') '7.03'
Set-TextValue $ws.Range('E11') '  +23.17%  '

Set-TextValue $ws.Range('D12') '0.603'
Set-TextValue $ws.Range('E12') '  -0.92%  '

Set-TextValue $ws.Range('D13') '48.12'
Set-TextValue $ws.Range('E13') '  -4.10%  '

Set-TextValue $ws.Range('D14') '0.0000282'
Set-TextValue $ws.Range('E14') '  -3.25%  '

Set-TextValue $ws.Range('D15') '4.217.43'
Set-TextValue $ws.Range('E15') '  +0.66%  '

Set-TextValue $ws.Range('D16') '673.38'
Set-TextValue $ws.Range('E16') '  -3.31%  '

Set-TextValue $ws.Range('D17') '8.86'
Set-TextValue $ws.Range('E17') '  -1.03%  '

Set-TextValue $ws.Range('D18') '3.634.54'
Set-TextValue $ws.Range('E18') '  +0.85%  '

Set-TextValue $ws.Range('D19') '70.571.50'
Set-TextValue $ws.Range('E19') '  -2.41%  '

Set-TextValue $ws.Range('E20') '  -0.23%  '

Set-TextValue $ws.Range('D21') '17.74'
Set-TextValue $ws.Range('E21') '  -4.28%  '

Set-TextValue $ws.Range('D22') '11.37'
Set-TextValue $ws.Range('E22') '  -3.14%  '

Set-TextValue $ws.Range('D23') '0.937'
Set-TextValue $ws.Range('E23') '  +0.32%  '

Set-TextValue $ws.Range('D24') '17.02'
Set-TextValue $ws.Range('E24') '  -3.85%  '

Set-TextValue $ws.Range('D25') '99.74'
Set-TextValue $ws.Range('E25') '  -5.04%  '

Set-TextValue $ws.Range('D26') '3.90'
Set-TextValue $ws.Range('E26') '  -3.32%  '

Set-TextValue $ws.Range('D27') '2.78'
Set-TextValue $ws.Range('E27') '  -2.13%  '

Set-TextValue $ws.Range('E28') '  +0.03%  '

Set-TextValue $ws.Range('D29') '9.83'
Set-TextValue $ws.Range('E29') '  -3.14%  '

Set-TextValue $ws.Range('D30') '34.46'
Set-TextValue $ws.Range('E30') '  -2.13%  '

Set-TextValue $ws.Range('B31') 'Filecoin'
Set-TextValue $ws.Range('C31') 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws.Range('D31') '8.97'
Set-TextValue $ws.Range('E31') '  -1.03%  '

Set-TextValue $ws.Range('B32') 'Stacks'
Set-TextValue $ws.Range('C32') 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue $ws.Range('D32') '3.30'
Set-TextValue $ws.Range('E32') '  -5.40%  '

Set-TextValue $ws.Range('D33') '7.46'
Set-TextValue $ws.Range('E33') '  -0.50%  '

Set-TextValue $ws.Range('E34') '  -7.65%  '

Set-TextValue $ws.Range('D35') '3.93'
Set-TextValue $ws.Range('E35') '  -5.48%  '

Set-TextValue $ws.Range('D36') '578.25'
Set-TextValue $ws.Range('E36') '  -2.96%  '

Set-TextValue $ws.Range('D37') '11.03'
Set-TextValue $ws.Range('E37') '  -2.76%  '

Set-TextValue $ws.Range('E38') '  -0.86%  '

Set-TextValue $ws.Range('D39') '58.16'
Set-TextValue $ws.Range('E39') '  -3.22%  '

Set-TextValue $ws.Range('E40') '  +0.05%  '

Set-TextValue $ws.Range('B41') 'Maker'
Set-TextValue $ws.Range('C41') 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue $ws.Range('D41') '3.549.16'
Set-TextValue $ws.Range('E41') '  -3.30%  '

Set-TextValue $ws.Range('B42') 'VeChain'
Set-TextValue $ws.Range('C42') 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range('D42') '0.0451'
Set-TextValue $ws.Range('E42') '  +2.36%  '

Set-TextValue $ws.Range('E43') '  -3.38%  '

Set-TextValue $ws.Range('D44') '0.343'
Set-TextValue $ws.Range('E44') '  -1.84%  '

Set-TextValue $ws.Range('D45') '34.42'
Set-TextValue $ws.Range('E45') '  -4.95%  '

Set-TextValue $ws.Range('D46') '0.0₃0729'
Set-TextValue $ws.Range('E46') '  -7.44%  '

Set-TextValue $ws.Range('D47') '2.67'
Set-TextValue $ws.Range('E47') '  -6.06%  '

Set-TextValue $ws.Range('D48') '2.83'
Set-TextValue $ws.Range('E48') '  +1.92%  '

Set-TextValue $ws.Range('E49') '  +0.20%  '

Set-TextValue $ws.Range('D50') '136.29'
Set-TextValue $ws.Range('E50') '  +1.83%  '

Set-TextValue $ws.Range('E51') '  -1.65%  '
